$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the two new columns (I and J)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the header style used by the other header cells (e.g. H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data values for columns I (I0) and J (IF), rows 2-29
$data = @{
    2  = @(1, 5)
    3  = @(1, 5)
    4  = @(1, 4)
    5  = @(1, 7)
    6  = @(1, 7)
    7  = @(1, 4)
    8  = @(1, 5)
    9  = @(1, 7)
    10 = @(1, 5)
    11 = @(1, 5)
    12 = @(1, 6)
    13 = @(1, 5)
    14 = @(1, 6)
    15 = @(1, 6)
    16 = @(1, 6)
    17 = @(1, 7)
    18 = @(1, 5)
    19 = @(1, 6)
    20 = @(1, 5)
    21 = @(1, 7)
    22 = @(1, 6)
    23 = @(1, 6)
    24 = @(1, 5)
    25 = @(1, 4)
    26 = @(5, 7)
    27 = @(1, 2)
    28 = @(3, 4)
    29 = @(3, 3)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
